$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (and values, temporarily) from column R (rows 3-14) into the
# new column S so the new cells inherit the same number formats, fonts, and
# borders as the corresponding row in column R.
$ws.Range("R3:R14").Copy($ws.Range("S3:S14"))

# Now overwrite with the real 2023 figures.
$ws.Range("S3").Value = 2023
$ws.Range("S4").Value = 78
$ws.Range("S5").Value = 77.400000000000006
$ws.Range("S6").Value = 1739
$ws.Range("S7").Value = 1631
$ws.Range("S8").Value = 1093.2
$ws.Range("S9").Value = 33.200000000000003
$ws.Range("S10").Value = 21.9
$ws.Range("S11").Value = 44.2
$ws.Range("S12").Value = 2.5
$ws.Range("S13").Value = 33.799999999999997
$ws.Range("S14").Value = 0.6

# Restore the original selection shown after the edit.
$ws.Range("G18").Select() | Out-Null
